# Append the "Mohammad" / Haptic Glove trial-notes block to the end of the
# document body, just before the trailing blank paragraph and the sectPr.
#
# Target shape (see diff):
#   ...Slight difference ... vibrations.</w:p>
#   <w:p/>                                                  (blank line)
#   <w:p>Mohammad</w:p>
#   <w:p>Angry and Surprise were easy to recognize if they
#         followed one another. If there was a gap, it was
#         confusing.</w:p>
#   <w:p/>                                                  (blank line)
#   <w:p/>                                                  (pre-existing trailing blank, untouched)
#   <w:sectPr>...

$d = $word.ActiveDocument
$wordNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# The last real (non-blank) paragraph before the document's final blank
# paragraph / sectPr is "Slight difference in top and middle ...".
$lastParaIndex = $d.Paragraphs.Count - 1
$anchor = $d.Paragraphs($lastParaIndex)

# 1) Blank paragraph right after the anchor paragraph.
$r = $anchor.Range
$r.Collapse(0)
$r.InsertParagraphAfter()

# 2) Another new (still-blank) paragraph that will become "Mohammad".
$r = $d.Paragraphs($lastParaIndex + 1).Range
$r.Collapse(0)
$r.InsertParagraphAfter()

# Make paragraph (lastParaIndex + 1) a plain empty <w:p/> (no stray run).
$blank1 = $d.Paragraphs($lastParaIndex + 1).Range
$null = $blank1.InsertXML('<w:p ' + $wordNs + '/>')

# 3) Fill the "Mohammad" paragraph.
$r = $d.Paragraphs($lastParaIndex + 2).Range
$r.Collapse(1)
$r.InsertAfter("Mohammad")

# 4) New blank paragraph after "Mohammad" that will become the Angry/Surprise text.
$r = $d.Paragraphs($lastParaIndex + 2).Range
$r.Collapse(0)
$r.InsertParagraphAfter()

# 5) Fill the "Angry and Surprise..." paragraph.
$r = $d.Paragraphs($lastParaIndex + 3).Range
$r.Collapse(1)
$r.InsertAfter("Angry and Surprise were easy to recognize if they followed one another. If there was a gap, it was confusing.")

# 6) Trailing blank paragraph after the new text block (still before the
#    document's original trailing blank paragraph).
$r = $d.Paragraphs($lastParaIndex + 3).Range
$r.Collapse(0)
$r.InsertParagraphAfter()

# Make that new trailing paragraph a plain empty <w:p/> (no stray run).
$blank2 = $d.Paragraphs($lastParaIndex + 4).Range
$null = $blank2.InsertXML('<w:p ' + $wordNs + '/>')
